$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for the new "Post Treatment" (column D) outcome measures, rows 2-15
$values = @(
    "A little worse ",
    "A little worse ",
    "Somewhat worse",
    "A lot worse",
    "A little worse ",
    "A little worse ",
    "Somewhat worse",
    "Somewhat worse",
    "A little worse ",
    "A little worse ",
    "A little worse ",
    "A little worse ",
    "A little worse ",
    "A lot worse"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Widen column D to fit the new content (stored width 19.5 in the XML)
$ws.Columns.Item(4).ColumnWidth = 18.666666666666668

# Update the active selection to reflect the next empty cell below the data
$ws.Range("D16").Select()
